# "clean data fixes my finals round numbering mistake."
#
# Column JU (the last populated column, one cell per row) holds the data
# for the final finals-round match. The round numbering for the finals
# was off, so the same match data needs to be duplicated across three
# more columns (JV, JW, JX) appended immediately after JU - one column
# per corrected finals match - while every row keeps the same per-row
# value it already had in JU.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 102
$srcCol  = 281   # JU - last existing column
$col1    = 282   # JV - new
$col2    = 283   # JW - new
$col3    = 284   # JX - new (becomes the new last column)

$srcRange = $ws.Range($ws.Cells.Item(1, $srcCol), $ws.Cells.Item($lastRow, $srcCol))
$values = $srcRange.Value()

# JU used to be the final column and carried no explicit cell style.
# Re-create its cells so they pick up the normal column formatting that
# every other data column already has (matching the rest of each row).
$srcRange.ClearContents()
$srcRange.Value = $values

$dst1 = $ws.Range($ws.Cells.Item(1, $col1), $ws.Cells.Item($lastRow, $col1))
$dst1.Value = $values

$dst2 = $ws.Range($ws.Cells.Item(1, $col2), $ws.Cells.Item($lastRow, $col2))
$dst2.Value = $values

# JX becomes the new final column, so (matching how JU looked before the
# fix) it is left without the normal column-level formatting.
$ws.Columns($col3).ClearFormats()
$dst3 = $ws.Range($ws.Cells.Item(1, $col3), $ws.Cells.Item($lastRow, $col3))
$dst3.Value = $values

Write-Host "Duplicated column JU into JV, JW, JX across rows 1-$lastRow"
